$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.495.15"
$ws.Range("E2").Value = "  -1.58%  "

$ws.Range("D3").Value = "1.910.54"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("D4").Value = "0.9990"

$ws.Range("D5").Value = "239.67"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "0.4744"
$ws.Range("E7").Value = "  -2.24%  "

$ws.Range("D8").Value = "0.2850"
$ws.Range("E8").Value = "  -3.18%  "

$ws.Range("D9").Value = "0.06688"
$ws.Range("E9").Value = "  -6.16%  "

$ws.Range("D10").Value = "18.80"
$ws.Range("E10").Value = "  -4.25%  "

$ws.Range("D11").Value = "100.87"
$ws.Range("E11").Value = "  -6.25%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.912.56"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07680"
$ws.Range("E13").Value = "  -1.21%  "

$ws.Range("D14").Value = "5.227"
$ws.Range("E14").Value = "  -2.88%  "

$ws.Range("D15").Value = "0.6704"
$ws.Range("E15").Value = "  -4.95%  "

$ws.Range("D16").Value = "30.480.74"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "256.83"
$ws.Range("E17").Value = "  -7.88%  "

$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").Value = "0.000007474"
$ws.Range("E19").Value = "  -4.77%  "

$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  -5.07%  "

$ws.Range("D21").Value = "5.397"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D23").Value = "0.4510"
$ws.Range("E23").Value = "  -8.83%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.332"
$ws.Range("E24").Value = "  -2.86%  "

$ws.Range("D25").Value = "168.47"
$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.379"
$ws.Range("E26").Value = "  -3.93%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.99"
$ws.Range("E27").Value = "  -3.72%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.056"
$ws.Range("E28").Value = "  -5.53%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "0.1011"
$ws.Range("E29").Value = "  -3.82%  "

$ws.Range("D30").Value = "1.376"
$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.669"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.513"
$ws.Range("E32").Value = "  -3.57%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.259"
$ws.Range("E33").Value = "  -3.65%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04722"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7304"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.114"
$ws.Range("E36").Value = "  -4.75%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "0.9988"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.704"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01922"
$ws.Range("E39").Value = "  -4.15%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.598"
$ws.Range("E40").Value = "  -2.99%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.271"
$ws.Range("E41").Value = "  -3.71%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "74.14"
$ws.Range("E42").Value = "  -4.79%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.965"
$ws.Range("E43").Value = "  -7.23%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8611"
$ws.Range("E44").Value = "  -3.93%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "105.94"
$ws.Range("E45").Value = "  -3.33%  "

$ws.Range("D46").Value = "0.4247"
$ws.Range("E46").Value = "  -4.67%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9990"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.000.77"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.418"
$ws.Range("E49").Value = "  -5.86%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1200"
$ws.Range("E50").Value = "  -4.10%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "34.78"
$ws.Range("E51").Value = "  -3.40%  "
